$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Controle de qualidade correto"
$ws.Range("B3").Value = "Controle de qualidade incorreto"
$ws.Range("B4").Value = "Relatório de controle de qualidade"
$ws.Range("B5").Value = "Baixar relatório"
$ws.Range("B6").Value = "O relatório de controle de qualidade não encontrou erros nos dados."
$ws.Range("B7").Value = "O relatório de controle de qualidade encontrou erros nos dados. Revise os detalhes no relatório."
$ws.Range("B8").Value = "Configuração correta"
$ws.Range("B9").Value = "Configuração incorreta"
$ws.Range("B10").Value = "Todos os pacotes foram instalados"
$ws.Range("B11").Value = "Alguns pacotes não foram instalados:"
$ws.Range("B12").Value = "Fechar"
